# "planned months now in Excel, and directory was reorganized"
#
# 1. Insert a new "Planned Months" worksheet between "Records" and
#    "School Exceptions" (it becomes sheetId 3, while "School Exceptions"
#    keeps sheetId 2 but moves to the 3rd tab / a new relationship id).
# 2. The new sheet lists the months the program runs, with a bold /
#    yellow-filled header cell (matching the existing "School Exceptions"
#    header style already used in this workbook).
# 3. The "Planned Months" tab becomes the active / selected tab.
# 4. The "Records" sheet's view scrolls over / selection moves (no longer
#    the selected tab) and loses its old selection-on-row-60.

$wb = $excel.ActiveWorkbook

$records = $wb.Worksheets.Item("Records")

# Make sure we're working from the Records sheet before touching its view,
# then move its selection / visible area.
$records.Activate()
$records.Range("E9").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# Insert the new "Planned Months" sheet right after "Records" (i.e. before
# "School Exceptions").
$newSheet = $wb.Worksheets.Add([Type]::Missing, $records)
$newSheet.Name = "Planned Months"

$months = @("February", "March", "April", "May", "June")

$header = $newSheet.Range("A1")
$header.Value = "Planned Months"
$header.Font.Bold = $true
$header.Interior.Color = 65535   # yellow, same header style as "School Exceptions"

for ($i = 0; $i -lt $months.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $months[$i]
}

$newSheet.Range("A7").Select()

# "Planned Months" becomes the active tab.
$newSheet.Activate()
